$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$vals = @(0, 0.1, 0.5, 0.9, 1, 1.01, 1.1, 1.5, 1.9, 1.99, 2, 2.5, 3, 3.99, 4, 8.43, 10, 10.5, 18, 18.5, 18.66, 18.666666666666668, 19, 19.33, 19.330729166666668, 19.5, 20, 25, 27, 100)
for ($i = 0; $i -lt $vals.Length; $i++) {
  $ws.Columns.Item($i+2).ColumnWidth = $vals[$i]
}
